$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style used by G1 onto the new H1 header cell
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
